# Daily refresh of the cryptos price table (Coin/Link/Price/Volume columns).
# Most rows keep the same coin but get updated Price (D) / Volume(1h) (E)
# values; rows 16-17 additionally swap which coin (TRON vs WrappedEther)
# occupies which row, along with its Link/Price/Volume.
#
# Several Price values look like plain numbers (e.g. "1.00", "27.71").
# Assigning such strings straight to .Value lets Excel auto-convert them to
# numeric cells, which silently rewrites the text (loses trailing zeros,
# introduces floating point noise, etc). Set-TextValue forces those specific
# cells to a text number format first so the original string is preserved
# exactly, matching the source inlineStr content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

$ws.Range('D2').Value = '65.037.15'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '3.524.06'
$ws.Range('E3').Value = '  -1.21%  '
$ws.Range('E4').Value = '  +0.03%  '
Set-TextValue 'D5' '594.82'
$ws.Range('E5').Value = '  -1.00%  '
Set-TextValue 'D6' '134.62'
$ws.Range('E6').Value = '  -2.43%  '
$ws.Range('D7').Value = '3.521.11'
$ws.Range('E7').Value = '  -1.26%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  -1.38%  '
$ws.Range('E10').Value = '  +0.72%  '
$ws.Range('E11').Value = '  +2.21%  '
Set-TextValue 'D12' '0.388'
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('D13').Value = '4.129.38'
$ws.Range('E13').Value = '  -1.09%  '
Set-TextValue 'D14' '27.71'
$ws.Range('E14').Value = '  +1.54%  '
$ws.Range('E15').Value = '  -0.64%  '
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D16' '0.117'
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.528.92'
$ws.Range('E17').Value = '  -1.09%  '
$ws.Range('D18').Value = '65.029.42'
$ws.Range('E18').Value = '  +0.10%  '
Set-TextValue 'D19' '10.11'
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('E20').Value = '  +0.09%  '
Set-TextValue 'D21' '5.72'
$ws.Range('E21').Value = '  -2.63%  '
Set-TextValue 'D22' '392.47'
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('E23').Value = '  +0.30%  '
$ws.Range('D24').Value = '3.669.66'
$ws.Range('E24').Value = '  -1.15%  '
Set-TextValue 'D25' '74.64'
$ws.Range('E25').Value = '  +0.67%  '
Set-TextValue 'D26' '0.999'
$ws.Range('E26').Value = '  -0.25%  '
Set-TextValue 'D28' '7.77'
$ws.Range('E28').Value = '  +0.74%  '
$ws.Range('E29').Value = '  +10.07%  '
Set-TextValue 'D30' '1.00'
$ws.Range('E30').Value = '  -0.15%  '
$ws.Range('E31').Value = '  -1.02%  '
$ws.Range('E32').Value = '  +0.64%  '
$ws.Range('D33').Value = '3.530.76'
$ws.Range('E33').Value = '  -1.29%  '
Set-TextValue 'D34' '24.21'
$ws.Range('E34').Value = '  +0.55%  '
$ws.Range('E35').Value = '  -0.01%  '
Set-TextValue 'D36' '0.144'
$ws.Range('E36').Value = '  -0.45%  '
Set-TextValue 'D37' '5.29'
$ws.Range('E37').Value = '  +5.31%  '
$ws.Range('E38').Value = '  +0.55%  '
Set-TextValue 'D39' '1.58'
$ws.Range('E39').Value = '  +1.90%  '
Set-TextValue 'D40' '168.72'
Set-TextValue 'D41' '0.0817'
$ws.Range('E41').Value = '  +0.86%  '
Set-TextValue 'D42' '0.825'
$ws.Range('E42').Value = '  -0.58%  '
$ws.Range('E43').Value = '  +4.01%  '
Set-TextValue 'D44' '25.97'
$ws.Range('E44').Value = '  -3.54%  '
Set-TextValue 'D45' '42.96'
$ws.Range('E45').Value = '  +0.70%  '
$ws.Range('E46').Value = '  +0.08%  '
Set-TextValue 'D47' '4.44'
$ws.Range('E47').Value = '  -0.72%  '
$ws.Range('E48').Value = '  -0.24%  '
$ws.Range('E49').Value = '  +0.19%  '
$ws.Range('D50').Value = '2.415.91'
$ws.Range('E50').Value = '  -2.21%  '
Set-TextValue 'D51' '0.911'
$ws.Range('E51').Value = '  +6.20%  '
